# "Generate Report for handoff"
#
# The localization-status report is regenerated: the source markdown file that is
# being tracked changed identity (new guid-named file), and the handoff for it
# failed (transform failure) instead of succeeding, so the per-language detail
# sheets ("zh-cn" and "de-de") no longer show a completed handoff: the "Latest
# Handoff File" link/column is cleared, the "Latest Handoff Datetime" reverts to
# the epoch placeholder, and the "Handoff Reason" flips from "Include" to
# "Ignored". The "Overview"/zh-cn/de-de "Status" column text changes from
# "Ready for handoff" to "Handoff transform failed", and every hyperlink that
# displayed the old markdown file name now displays the new one.

$wb = $excel.ActiveWorkbook

$oldFileName = "b4cd3aed-69e7-4617-a156-447920c7b6c3.md"
$newFileName = "d371a348-48f4-4032-9075-53a315364416.md"

$oldStatus = "Ready for handoff"
$newStatus = "Handoff transform failed"

$epoch = "0001-01-01 00:00:00"
$ignored = "Ignored"

# 1) Update the status text and the tracked source file name (both the cell
#    text and the hyperlink display text) on every sheet. The hyperlink target
#    addresses themselves are untouched - only the visible text / display
#    changes.
foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($cell in $ws.UsedRange.Cells) {
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        } elseif ($cell.Value2 -eq $oldFileName) {
            $cell.Value = $newFileName
        }
    }

    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.TextToDisplay -eq $oldFileName) {
            $hl.TextToDisplay = $newFileName
        }
    }
}

# 2) The per-language detail sheets ("zh-cn", "de-de") recorded a completed
#    handoff (link to the generated .xlf, a real handoff datetime, and reason
#    "Include"). Since the transform failed this run, that handoff info is
#    removed / reset on row 2 (the tracked source file row).
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Drop the "Latest Handoff File" hyperlink + cell contents in C2.
    $linksToDelete = @()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$C$2') {
            $linksToDelete += $hl
        }
    }
    foreach ($hl in $linksToDelete) {
        $hl.Delete()
    }
    $ws.Range("C2").Clear()

    # "Latest Handoff Datetime" resets to the epoch placeholder.
    $ws.Range("D2").Value = $epoch

    # "Handoff Reason" flips from Include to Ignored.
    $ws.Range("H2").Value = $ignored
}
